$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new text labels first, in the order they appear in the shared string
# table of the target workbook (row 13, 12, 11, 8, 9/10).
$ws.Range("A13").Value = "Работа по созданию функционала загрузки данных (DatePicker, Валидация Int, Float, date)"
$ws.Range("A12").Value = "Работа по созданию функционала загрузки данных (Разработка строки ввода с выбором элемента сравнения)"
$ws.Range("A11").Value = "Доработка Side-Button открывающей фильтры"
$ws.Range("A8").Value = "Работа по созданию функционала загрузки данных (Разработка каркаса приложения, структуры филтров)"
$ws.Range("A9").Value = "Работа по созданию функционала загрузки данных (Разработка, концепции, внешнего вида)"
$ws.Range("A10").Value = "Работа по созданию функционала загрузки данных (Разработка, концепции, внешнего вида)"

# Row 7: was "В паспортах указать ДОли компании в общих закупках (продажах)" with no time/date.
# Now becomes the "Анализ, формированеи рынков для организации ЛИГА-7" entry (used to be row 8).
$ws.Range("A7").Value = "Анализ, формированеи рынков для организации ЛИГА-7"
$ws.Range("B7").Value = 1
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 43512

# Row 8: "Работа по созданию функционала загрузки данных (Разработка каркаса приложения, структуры филтров)"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 43510

# Row 9: "Работа по созданию функционала загрузки данных (Разработка, концепции, внешнего вида)"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 43512

# Row 10: same text as row 9, B/C updated
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 43513

# Row 11: "Доработка Side-Button открывающей фильтры"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 43514

# Row 12: previously blank, now filled in with date style copied from C4
$ws.Range("B12").Value = 4
$ws.Range("C4").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = 43516

# Row 13: previously blank, now filled in with date style copied from C4
$ws.Range("B13").Value = 2
$ws.Range("C4").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = 43517

# Update the sheet view: select A10 instead of D39, drop the frozen topLeftCell="A3"
$ws.Range("A10").Select()
